$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1762.65
$ws.Range("J17").Value = 1762.65
$ws.Range("L17").Value = 5287.950000000001
$ws.Range("N17").Value = -5623.950000000001
# Row 19
$ws.Range("H19").Value = 847.5
$ws.Range("I19").Value = 979.25
$ws.Range("J19").Value = 689.4
$ws.Range("K19").Value = 979.25
$ws.Range("L19").Value = 689.4
$ws.Range("M19").Value = -804.25
$ws.Range("N19").Value = -1039.4
# Row 98
$ws.Range("H98").Value = 642.7692
$ws.Range("I98").Value = 642.7692
$ws.Range("K98").Value = 642.7692
$ws.Range("M98").Value = 855.2308
# Row 112
$ws.Range("H112").Value = 2022.7878
$ws.Range("I112").Value = 1254.8572
$ws.Range("J112").Value = 2229.5386
$ws.Range("K112").Value = 3764.5716
$ws.Range("L112").Value = 6688.6158
$ws.Range("M112").Value = -2656.5716
$ws.Range("N112").Value = -8904.6158
# Row 122
$ws.Range("H122").Value = 642.7692
$ws.Range("I122").Value = 642.7692
$ws.Range("K122").Value = 1928.3076
$ws.Range("M122").Value = 521.6924000000001
# Row 134
$ws.Range("H134").Value = 98333.336
$ws.Range("J134").Value = 98333.336
$ws.Range("L134").Value = 98333.336
$ws.Range("N134").Value = -108473.336

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1005.6125
$ws.Range("I2").Value = 819.71014
$ws.Range("J2").Value = 2171.7273
$ws.Range("K2").Value = 819.71014
$ws.Range("L2").Value = 2171.7273
$ws.Range("M2").Value = -706.71014
$ws.Range("N2").Value = -2397.7273
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").Value = $null
# Row 61
$ws.Range("H61").Value = 4379.6665
$ws.Range("I61").Value = 4337.375
$ws.Range("J61").Value = 4464.25
$ws.Range("K61").Value = 4337.375
$ws.Range("L61").Value = 4464.25
$ws.Range("M61").Value = -4125.375
$ws.Range("N61").Value = -4888.25
# Row 116
$ws.Range("H116").Value = 1005.6125
$ws.Range("I116").Value = 819.71014
$ws.Range("J116").Value = 2171.7273
$ws.Range("K116").Value = 819.71014
$ws.Range("L116").Value = 2171.7273
$ws.Range("M116").Value = 1474.28986
$ws.Range("N116").Value = -6759.7273
# Row 122
$ws.Range("H122").Value = 2112.5366
$ws.Range("I122").Value = 1960.75
$ws.Range("K122").Value = 5882.25
$ws.Range("M122").Value = -3432.25
# Row 132
$ws.Range("H132").Value = 2164.8572
$ws.Range("I132").Value = 1882.6
$ws.Range("J132").Value = 2870.5
$ws.Range("K132").Value = 5647.799999999999
$ws.Range("L132").Value = 8611.5
$ws.Range("M132").Value = -3117.799999999999
$ws.Range("N132").Value = -13671.5
# Row 136
$ws.Range("H136").Value = 4379.6665
$ws.Range("I136").Value = 4337.375
$ws.Range("J136").Value = 4464.25
$ws.Range("K136").Value = 13012.125
$ws.Range("L136").Value = 13392.75
$ws.Range("M136").Value = -10462.125
$ws.Range("N136").Value = -18492.75

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1005.6125
$ws.Range("I3").Value = 819.71014
$ws.Range("J3").Value = 2171.7273
$ws.Range("K3").Value = 819.71014
$ws.Range("L3").Value = 2171.7273
$ws.Range("M3").Value = -705.71014
$ws.Range("N3").Value = -2399.7273
# Row 99
$ws.Range("H99").Value = 1775
$ws.Range("I99").Value = 1775
$ws.Range("K99").Value = 1775
$ws.Range("M99").Value = -277
# Row 134
$ws.Range("H134").Value = 6170.5884
$ws.Range("I134").Value = 4191.6
$ws.Range("J134").Value = 8997.714
$ws.Range("K134").Value = 12574.8
$ws.Range("L134").Value = 26993.142
$ws.Range("M134").Value = -10039.8
$ws.Range("N134").Value = -32063.142

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 416.375
$ws.Range("I22").Value = 238.16667
$ws.Range("K22").Value = 238.16667
$ws.Range("M22").Value = 111.83333
# Row 31
$ws.Range("H31").Value = 5893.391
$ws.Range("I31").Value = 4549.8
$ws.Range("K31").Value = 4549.8
$ws.Range("M31").Value = -4254.8
# Row 34
$ws.Range("H34").Value = 5893.391
$ws.Range("I34").Value = 4549.8
$ws.Range("K34").Value = 4549.8
$ws.Range("M34").Value = -4347.8
# Row 99
$ws.Range("H99").Value = 2564.8333
$ws.Range("I99").Value = 2564.8333
$ws.Range("K99").Value = 2564.8333
$ws.Range("M99").Value = -1066.8333
# Row 108
$ws.Range("H108").Value = 57500
$ws.Range("I108").Value = 75000
$ws.Range("J108").Value = 40000
$ws.Range("K108").Value = 75000
$ws.Range("L108").Value = 40000
$ws.Range("M108").Value = -71160
$ws.Range("N108").Value = -47680
# Row 122
$ws.Range("H122").Value = 4427.857
$ws.Range("I122").Value = 4427.857
$ws.Range("K122").Value = 13283.571
$ws.Range("M122").Value = -10833.571
# Row 126
$ws.Range("H126").Value = 2564.8333
$ws.Range("I126").Value = 2564.8333
$ws.Range("K126").Value = 7694.499899999999
$ws.Range("M126").Value = -5224.499899999999

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 60345.74
$ws.Range("J131").Value = 2492.5454
$ws.Range("L131").Value = 7477.6362
$ws.Range("N131").Value = -17557.6362
# Row 137
$ws.Range("H137").Value = 5736.7334
$ws.Range("J137").Value = 2594.889
$ws.Range("L137").Value = 7784.667
$ws.Range("N137").Value = -17984.667

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 224.5
$ws.Range("J2").Value = 149
$ws.Range("L2").Value = 149
$ws.Range("N2").Value = -375
# Row 63
$ws.Range("H63").Value = 38999.5
$ws.Range("I63").Value = 38999.5
$ws.Range("K63").Value = 38999.5
$ws.Range("M63").Value = -38313.5
# Row 66
$ws.Range("H66").Value = 38999.5
$ws.Range("I66").Value = 38999.5
$ws.Range("K66").Value = 116998.5
$ws.Range("M66").Value = -113566.5
# Row 102
$ws.Range("H102").Value = 26038.71
$ws.Range("I102").Value = 28544.908
$ws.Range("K102").Value = 28544.908
$ws.Range("M102").Value = -26922.908
# Row 122
$ws.Range("H122").Value = 64706.25
$ws.Range("I122").Value = 126612.625
$ws.Range("K122").Value = 379837.875
$ws.Range("M122").Value = -377387.875
# Row 126
$ws.Range("H126").Value = 1999
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = $null
# Row 132
$ws.Range("H132").Value = 3334
$ws.Range("I132").Value = 4272.25
$ws.Range("K132").Value = 12816.75
$ws.Range("M132").Value = -10286.75

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 7798.5713
$ws.Range("I7").Value = 7129.1
$ws.Range("J7").Value = 9472.25
$ws.Range("K7").Value = 7129.1
$ws.Range("L7").Value = 9472.25
$ws.Range("M7").Value = -7017.1
$ws.Range("N7").Value = -9696.25
# Row 20
$ws.Range("H20").Value = 9190.9375
$ws.Range("I20").Value = 7000
$ws.Range("K20").Value = 7000
$ws.Range("M20").Value = -6774
# Row 40
$ws.Range("H40").Value = 3479.4348
$ws.Range("I40").Value = 2944
$ws.Range("K40").Value = 2944
$ws.Range("M40").Value = -2808
# Row 61
$ws.Range("H61").Value = 1162.7059
$ws.Range("I61").Value = 922.875
$ws.Range("K61").Value = 922.875
$ws.Range("M61").Value = -720.875
# Row 68
$ws.Range("H68").Value = 2837.6924
$ws.Range("I68").Value = 2444.5454
$ws.Range("K68").Value = 2444.5454
$ws.Range("M68").Value = -1695.5454
# Row 71
$ws.Range("H71").Value = 2837.6924
$ws.Range("I71").Value = 2444.5454
$ws.Range("K71").Value = 12222.727
$ws.Range("M71").Value = -8478.726999999999
# Row 93
$ws.Range("H93").Value = 3556.5
$ws.Range("I93").Value = 3857.1
$ws.Range("J93").Value = 2805
$ws.Range("K93").Value = 3857.1
$ws.Range("L93").Value = 2805
$ws.Range("M93").Value = -2609.1
$ws.Range("N93").Value = -5301
# Row 113
$ws.Range("H113").Value = 1162.7059
$ws.Range("I113").Value = 922.875
$ws.Range("K113").Value = 922.875
$ws.Range("M113").Value = 1247.125
# Row 122
$ws.Range("H122").Value = 3766.3333
$ws.Range("J122").Value = 4583.3335
$ws.Range("L122").Value = 13750.0005
$ws.Range("N122").Value = -18650.0005
# Row 126
$ws.Range("H126").Value = 7798.5713
$ws.Range("I126").Value = 7129.1
$ws.Range("J126").Value = 9472.25
$ws.Range("K126").Value = 21387.3
$ws.Range("L126").Value = 28416.75
$ws.Range("M126").Value = -18917.3
$ws.Range("N126").Value = -33356.75
# Row 132
$ws.Range("H132").Value = 2467.2083
$ws.Range("I132").Value = 2548.1052
$ws.Range("K132").Value = 7644.3156
$ws.Range("M132").Value = -5114.3156

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 751.44446
$ws.Range("I107").Value = 456
$ws.Range("J107").Value = 987.8
$ws.Range("K107").Value = 1368
$ws.Range("L107").Value = 2963.4
$ws.Range("M107").Value = 552
$ws.Range("N107").Value = -6803.4
# Row 122
$ws.Range("H122").Value = 14016.474
$ws.Range("I122").Value = 16371.077
$ws.Range("J122").Value = 8914.833000000001
$ws.Range("K122").Value = 49113.231
$ws.Range("L122").Value = 26744.499
$ws.Range("M122").Value = -46663.231
$ws.Range("N122").Value = -31644.499
# Row 126
$ws.Range("H126").Value = 2454.5
$ws.Range("I126").Value = 2394.111
$ws.Range("K126").Value = 7182.333
$ws.Range("M126").Value = -4712.333
# Row 132
$ws.Range("H132").Value = 2802
$ws.Range("I132").Value = 2602
$ws.Range("J132").Value = 3002
$ws.Range("K132").Value = 7806
$ws.Range("L132").Value = 9006
$ws.Range("M132").Value = -5276
$ws.Range("N132").Value = -14066
# Row 136
$ws.Range("H136").Value = 2180.6155
$ws.Range("I136").Value = 2516.4443
$ws.Range("J136").Value = 1425
$ws.Range("K136").Value = 7549.3329
$ws.Range("L136").Value = 4275
$ws.Range("M136").Value = -4999.3329
$ws.Range("N136").Value = -9375
